$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 44 (weekly update: newest "Poroto granado" price
# entry). Everything that was rows 44-68 shifts down to 45-69.
$ws.Rows.Item(44).EntireRow.Insert()
$ws.Range("A44").Value = 7
$ws.Range("B44").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C44").Value = "Ñuble"
$ws.Range("D44").Value = 44567
$ws.Range("E44").Value = 16
$ws.Range("F44").Value = 100112030
$ws.Range("G44").Value = "Poroto granado"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 120
$ws.Range("K44").Value = 30000
$ws.Range("L44").Value = 31000
$ws.Range("M44").Value = 30500
$ws.Range("N44").Value = "`$/saco 25 kilos"
$ws.Range("O44").Value = "Provincia de Diguillín"
$ws.Range("P44").Value = 1220
$ws.Range("Q44").Value = 25
$ws.Range("R44").Value = "Hortaliza"

# Insert a second new data row at row 68 (another weekly entry). Rows that are
# now 68-69 (originally 67-68 before this pass) shift down to 69-70.
$ws.Rows.Item(68).EntireRow.Insert()
$ws.Range("A68").Value = 7
$ws.Range("B68").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C68").Value = "Ñuble"
$ws.Range("D68").Value = 44568
$ws.Range("E68").Value = 16
$ws.Range("F68").Value = 100112030
$ws.Range("G68").Value = "Poroto granado"
$ws.Range("H68").Value = "Sin especificar"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 120
$ws.Range("K68").Value = 30000
$ws.Range("L68").Value = 31000
$ws.Range("M68").Value = 30500
$ws.Range("N68").Value = "`$/saco 25 kilos"
$ws.Range("O68").Value = "Provincia de Diguillín"
$ws.Range("P68").Value = 1220
$ws.Range("Q68").Value = 25
$ws.Range("R68").Value = "Hortaliza"
